$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells we update stay stored as text, matching the original inlineStr cells
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.161.28'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.22'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.53'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5236'
$ws.Range("E7").Value = '  +1.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3773'
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07251'
$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  +0.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9035'
$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08567'
$ws.Range("E12").Value = '  +11.96%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.909.42'
$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.69'
$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.299'
$ws.Range("E15").Value = '  +0.52%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008639'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.56'
$ws.Range("E18").Value = '  +0.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.193.66'
$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.076'
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.152.67'
$ws.Range("E22").Value = '  +1.92%  '

$ws.Range("E23").Value = '  +0.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.441'
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("E25").Value = '  +1.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.16'
$ws.Range("E26").Value = '  +0.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.27'
$ws.Range("E27").Value = '  +1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.751'
$ws.Range("E28").Value = '  -0.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.09'
$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.826'
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.923'
$ws.Range("E31").Value = '  -0.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09305'
$ws.Range("E32").Value = '  +1.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8067'
$ws.Range("E33").Value = '  +2.81%  '

$ws.Range("E34").Value = '  -0.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.245'
$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.448'
$ws.Range("E36").Value = '  +4.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.951'
$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.618'
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5729'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01999'
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.140'
$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.648'
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '116.33'
$ws.Range("E44").Value = '  -1.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1519'
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.21'
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.620'
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.63'
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.28'
$ws.Range("E51").Value = '  +0.52%  '
